$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "GLOBAL RESULTS": refreshed computed values in column C
# ---------------------------------------------------------------------------
$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")
$wsGlobal.Cells.Item(2, 3).Value  = 45.98540759259373
$wsGlobal.Cells.Item(3, 3).Value  = 12.243126778022365
$wsGlobal.Cells.Item(4, 3).Value  = 0.7144207146200321
$wsGlobal.Cells.Item(6, 3).Value  = 26.286923460279134
$wsGlobal.Cells.Item(7, 3).Value  = 11.789585775420619
$wsGlobal.Cells.Item(8, 3).Value  = 0.7690953368472477
$wsGlobal.Cells.Item(10, 3).Value = 26.286923460279134
$wsGlobal.Cells.Item(11, 3).Value = 11.789585775420619
$wsGlobal.Cells.Item(12, 3).Value = 0.7690953368472477
$wsGlobal.Cells.Item(14, 3).Value = 30.337908673064945
$wsGlobal.Cells.Item(15, 3).Value = 11.882856297283656
$wsGlobal.Cells.Item(16, 3).Value = 0.4950125969433359
$wsGlobal.Cells.Item(18, 3).Value = 32.013360757965835
$wsGlobal.Cells.Item(19, 3).Value = 11.921432170091581
$wsGlobal.Cells.Item(20, 3).Value = 0.7143833406573219
$wsGlobal.Cells.Item(22, 3).Value = 12.193982195394145
$wsGlobal.Cells.Item(23, 3).Value = 42.87406608391055

# ---------------------------------------------------------------------------
# Sheet "WING": the Ycg estimation comparison block lost its erroneous
# "SFORZA" rows (rows 11 and 15) - the surviving "TORENBEEK_1982" rows move
# up, with row 11 now holding the former row 12 values and a new row 14
# holding the former row 16 values.
# ---------------------------------------------------------------------------
$wsWing = $wb.Worksheets.Item("WING")

# Drop the last two rows of the block (old rows 15 & 16, the second
# SFORZA/TORENBEEK_1982 pair) - nothing below them to shift.
$wsWing.Rows("15:16").Delete()

# Row 11 (was SFORZA / -0.131...) becomes TORENBEEK_1982 / 1.216893150576782
$wsWing.Cells.Item(11, 1).Value = "TORENBEEK_1982"
$wsWing.Cells.Item(11, 3).Value = 1.216893150576782

# Row 12 collapses to just the blank-separator string (old B12/C12 removed)
$wsWing.Cells.Item(12, 1).Value = " "
$wsWing.Cells.Item(12, 2).ClearContents()
$wsWing.Cells.Item(12, 3).ClearContents()

# Row 13 becomes the "Ycg ESTIMATION METHOD COMPARISON" header
$wsWing.Cells.Item(13, 1).Value = "Ycg ESTIMATION METHOD COMPARISON"

# Row 14 is rebuilt as the lone TORENBEEK_1982 result row
$wsWing.Cells.Item(14, 1).Value = "TORENBEEK_1982"
$wsWing.Cells.Item(14, 2).Value = "m"
$wsWing.Cells.Item(14, 3).Value = 4.735499999999998

# ---------------------------------------------------------------------------
# Sheet "VERTICAL TAIL": refreshed Ycg BRF (semi-tail) value
# ---------------------------------------------------------------------------
$wsVTail = $wb.Worksheets.Item("VERTICAL TAIL")
$wsVTail.Cells.Item(8, 3).Value = 3.1315999999999997

# ---------------------------------------------------------------------------
# Sheet "LANDING GEARS": refreshed Xcg BRF value
# ---------------------------------------------------------------------------
$wsGear = $wb.Worksheets.Item("LANDING GEARS")
$wsGear.Cells.Item(2, 3).Value = 12.297720537759226
